# Auto-generated edit script: update Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.283.06'
$ws.Range('E2').Value = '  -0.65%  '
$ws.Range('D3').Value = '3.094.53'
$ws.Range('E3').Value = '  +2.51%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '388.92'
$ws.Range('E5').Value = '  +2.35%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '103.65'
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.583'
$ws.Range('E9').Value = '  -1.14%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '36.91'
$ws.Range('E10').Value = '  +0.54%  '
$ws.Range('E11').Value = '  +0.15%  '
$ws.Range('E12').Value = '  -0.41%  '
$ws.Range('D13').Value = '3.582.62'
$ws.Range('E13').Value = '  +2.29%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '18.47'
$ws.Range('E14').Value = '  +0.22%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.73'
$ws.Range('E15').Value = '  +0.41%  '
$ws.Range('D16').Value = '3.098.52'
$ws.Range('E16').Value = '  +2.36%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.996'
$ws.Range('E17').Value = '  +2.28%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '10.75'
$ws.Range('E18').Value = '  +2.27%  '
$ws.Range('D19').Value = '51.373.82'
$ws.Range('E19').Value = '  -0.46%  '
$ws.Range('E20').Value = '  +5.16%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '12.43'
$ws.Range('E21').Value = '  +0.23%  '
$ws.Range('D22').Value = '0.0₃0964'
$ws.Range('E22').Value = '  +0.27%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '70.01'
$ws.Range('E23').Value = '  +0.09%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '265.59'
$ws.Range('E24').Value = '  -0.36%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.18'
$ws.Range('E25').Value = '  +0.99%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '7.96'
$ws.Range('E26').Value = '  -2.89%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '27.27'
$ws.Range('E27').Value = '  +4.19%  '
$ws.Range('E28').Value = '  -5.76%  '
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('E30').Value = '  -3.70%  '
$ws.Range('E31').Value = '  -1.60%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '10.47'
$ws.Range('E32').Value = '  +2.14%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '36.09'
$ws.Range('E33').Value = '  +7.11%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0473'
$ws.Range('E34').Value = '  +5.90%  '
$ws.Range('E35').Value = '  +0.10%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '49.73'
$ws.Range('E36').Value = '  -1.36%  '
$ws.Range('E37').Value = '  -0.09%  '
$ws.Range('E38').Value = '  +3.07%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.291'
$ws.Range('E39').Value = '  -2.57%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '130.49'
$ws.Range('E40').Value = '  +2.75%  '
$ws.Range('E41').Value = '  -0.15%  '
$ws.Range('E42').Value = '  +3.78%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '16.63'
$ws.Range('E43').Value = '  -1.97%  '
$ws.Range('E44').Value = '  +0.06%  '
$ws.Range('E45').Value = '  -1.85%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '22.09'
$ws.Range('E46').Value = '  +2.66%  '
$ws.Range('E47').Value = '  +5.01%  '
$ws.Range('E48').Value = '  -0.60%  '
$ws.Range('D49').Value = '2.070.59'
$ws.Range('E49').Value = '  +2.22%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.935'
$ws.Range('E50').Value = '  +19.06%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0325'
$ws.Range('E51').Value = '  +1.55%  '
